# BOM_LTC3556_v1.2a_.xlsx -- "Add files via upload" edit
#
# - Added TC7660EOA charge pump for AFE -2.5V AVSS:
#     Row 7 (USB connector line) gets its MFR part number filled in (col D).
# - Updated integration of both U1 and U2 IC's, V_Out I/O connections:
#     Row 21 (SPDT switch line) gets Manufacturer (col C) and MFR part
#     number (col D) filled in.
# - A handful of cosmetic column-width tweaks and the leftover "bold" cell
#   style on F1 is cleared back to Normal.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New / updated cell content -------------------------------------------

# Row 7: USB_B_Micro connector -- fill in the MFR Part number (col D) with
# the same footprint/part reference already used in col H.
$ws.Range("D7").Value = "Connector_USB:USB_Micro-B_Amphenol_10118194-0001LF_Horizontal"

# Row 21: SW_SPDT switch -- fill in Manufacturer (col C) and MFR Part
# number (col D).
$ws.Range("C21").Value = "C&K"
$ws.Range("D21").Value = "Button_Switch_THT:SW_Slide_SPDT_Straight_CK_OS102011MS2Q"

# --- Formatting -------------------------------------------------------------

# F1 had a vestigial "apply font" style left over from an earlier bold/unbold
# edit; clear it back to the Normal style.
$ws.Range("F1").Style = "Normal"

# Column width tweaks (values converted to the ColumnWidth property, which
# Excel stores internally with a small fixed padding offset).
$ws.Columns.Item(1).ColumnWidth = 16.764322916666668   # A
$ws.Columns.Item(2).ColumnWidth = 42.432291666666664   # B
$ws.Columns.Item(3).ColumnWidth = 42.432291666666664   # C
$ws.Columns.Item(4).ColumnWidth = 62.432291666666664   # D (widened)
$ws.Columns.Item(5).ColumnWidth = 27.565104166666668   # E
$ws.Columns.Item(6).ColumnWidth = 42.432291666666664   # F
$ws.Columns.Item(7).ColumnWidth = 81.76432291666667    # G
$ws.Columns.Item(8).ColumnWidth = 70.29947916666667    # H

# --- Selection ---------------------------------------------------------------
$ws.Range("D27").Select() | Out-Null
